# Natmi following Dr Hou advice
# Update Ligand-Receptor pair statistics for Bsg-Sele (rows 2-4) to reflect
# the revised ligand/receptor-expressing cell counts (1 -> 3) and the
# corresponding recalculated expression / specificity values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 49.400308
$ws.Range("H2").Value = 148.200924
$ws.Range("I2").Value = 0.3028101582105581
$ws.Range("J2").Value = 0.3028101582105581
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 7.321929333333333
$ws.Range("N2").Value = 21.965788
$ws.Range("Q2").Value = 361.7055642209013
$ws.Range("R2").Value = 3255.350077988111
$ws.Range("S2").Value = 0.3028101582105581
$ws.Range("T2").Value = 0.3028101582105581

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 69.564149
$ws.Range("H3").Value = 208.692447
$ws.Range("I3").Value = 0.4264088994034781
$ws.Range("J3").Value = 0.4264088994034782
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.321929333333333
$ws.Range("N3").Value = 21.965788
$ws.Range("Q3").Value = 509.3437831114707
$ws.Range("R3").Value = 4584.094048003236
$ws.Range("S3").Value = 0.4264088994034781
$ws.Range("T3").Value = 0.4264088994034782

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 44.17507666666666
$ws.Range("H4").Value = 132.52523
$ws.Range("I4").Value = 0.2707809423859638
$ws.Range("J4").Value = 0.2707809423859638
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 7.321929333333333
$ws.Range("N4").Value = 21.965788
$ws.Range("Q4").Value = 323.4467896479155
$ws.Range("R4").Value = 2911.02110683124
$ws.Range("S4").Value = 0.2707809423859638
$ws.Range("T4").Value = 0.2707809423859638
